$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Insert a new row at row 98 (pushes existing row 98+ down to 99+),
# so the new QLSECTIONPP command lands directly below the existing
# QLSECTION command (row 97).
$ws.Rows.Item(98).Insert()

# Update the description of the existing QLSECTION command (row 97)
# to reference the RFLTools DLL instead of the old theswamp.org link.
$ws.Range("B97").Value = "Draws a cross section from a Lidar point database (Note:  requires RFLTools DLL)"

# Populate the newly inserted row with the new QLSECTIONPP command.
$ws.Range("A98").Value = "QLSECTIONPP"
$ws.Range("B98").Value = "Draws a point to point section from a Lidar point database (Note:  requires RFLTools DLL)"

# Restore view state to match the committed workbook.
$ws.Application.ActiveWindow.ScrollRow = 47
$ws.Range("B98").Select()
